$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> Price (column D), Volume(1h) (column E) updated values.
# Only rows/columns that actually change per the diff are included;
# rows where the price cell is unchanged are left as $null for D.
$updates = @(
    @{ Row = 2;  D = "26.132.70";     E = "  -6.63%  " },
    @{ Row = 3;  D = "1.672.40";      E = "  -4.15%  " },
    @{ Row = 4;  D = "1.006";         E = "  +0.34%  " },
    @{ Row = 5;  D = "218.48";        E = "  -3.21%  " },
    @{ Row = 6;  D = "0.5072";        E = "  -12.54%  " },
    @{ Row = 7;  D = "1.007";         E = "  +0.24%  " },
    @{ Row = 8;  D = "0.2628";        E = "  -2.77%  " },
    @{ Row = 9;  D = "0.06325";       E = "  -3.55%  " },
    @{ Row = 10; D = "21.43";         E = "  -7.07%  " },
    @{ Row = 11; D = "0.07372";       E = "  -1.80%  " },
    @{ Row = 12; D = "1.673.66";      E = "  -4.26%  " },
    @{ Row = 13; D = "4.540";         E = "  -3.41%  " },
    @{ Row = 14; D = "0.5753";        E = "  -4.45%  " },
    @{ Row = 15; D = "1.899.42";      E = "  -4.15%  " },
    @{ Row = 16; D = "0.000008471";   E = "  -1.71%  " },
    @{ Row = 17; D = "64.72";         E = "  -12.31%  " },
    @{ Row = 18; D = "26.227.21";     E = "  -6.34%  " },
    @{ Row = 19; D = "4.957";         E = "  -6.64%  " },
    @{ Row = 20; D = "1.006";         E = "  +0.31%  " },
    @{ Row = 21; D = $null;           E = "  -3.90%  " },
    @{ Row = 22; D = "187.07";        E = "  -8.21%  " },
    @{ Row = 23; D = "6.176";         E = "  -6.79%  " },
    @{ Row = 24; D = "1.007";         E = "  +0.27%  " },
    @{ Row = 25; D = "143.06";        E = "  -4.87%  " },
    @{ Row = 26; D = "7.610";         E = "  -5.28%  " },
    @{ Row = 27; D = $null;           E = "  -5.09%  " },
    @{ Row = 28; D = "15.69";         E = "  -2.26%  " },
    @{ Row = 29; D = $null;           E = "  -5.47%  " },
    @{ Row = 30; D = "0.05761";       E = "  -5.18%  " },
    @{ Row = 31; D = "1.325";         E = "  -4.38%  " },
    @{ Row = 32; D = "3.501";         E = "  -6.08%  " },
    @{ Row = 33; D = $null;           E = "  -5.63%  " },
    @{ Row = 34; D = "1.663";         E = "  -0.40%  " },
    @{ Row = 35; D = "1.005";         E = "  -2.57%  " },
    @{ Row = 36; D = "0.5977";        E = "  -5.58%  " },
    @{ Row = 37; D = $null;           E = "  -3.42%  " },
    @{ Row = 38; D = "2.640";         E = "  -2.64%  " },
    @{ Row = 39; D = "0.01603";       E = "  -3.99%  " },
    @{ Row = 40; D = "1.080.47";      E = "  -3.81%  " },
    @{ Row = 41; D = "5.898";         E = "  -5.92%  " },
    @{ Row = 42; D = "0.8599";        E = "  -0.60%  " },
    @{ Row = 43; D = $null;           E = "  +0.00%  " },
    @{ Row = 44; D = "99.76";         E = "  +0.47%  " },
    @{ Row = 45; D = "1.821.45";      E = "  -3.89%  " },
    @{ Row = 46; D = $null;           E = "  +4.22%  " },
    @{ Row = 47; D = "56.05";         E = "  -5.10%  " },
    @{ Row = 48; D = "1.007";         E = "  +0.83%  " },
    @{ Row = 49; D = "8.068";         E = "  -2.45%  " },
    @{ Row = 50; D = "0.4306";        E = "  -2.78%  " },
    @{ Row = 51; D = "0.05200";       E = "  -3.53%  " }
)

# Rows whose new Price text would otherwise be re-interpreted by Excel's
# automatic number parsing (stripping the significant trailing zero, e.g.
# "4.540" -> 4.54). Force those specific cells to keep their literal text.
$forceTextRows = @(13, 16, 26, 38, 51)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        if ($forceTextRows -contains $u.Row) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
